# Realestate Update resale numbers 2025-02-02 20:18
# Append a new data row (row 48) to the CityResaleNum sheet, mirroring the
# layout of the existing rows (Date/Time/Weekday/Week as text, the 16 city
# columns as numbers, with -1 marking "no data").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 48

# --- Text columns (A-D) -----------------------------------------------
# A plain `.Value = "2025-02-02"` assignment gets auto-converted by Excel's
# usual "typed input" parsing into a date/number (same as typing it into the
# grid). Routing the literal through a TEXT() formula and then converting
# the formula to its computed value via Copy + PasteSpecial(values) keeps it
# as plain text without touching any cell's NumberFormat/style.
$ws.Range("A$row").Formula = '=TEXT("2025-02-02","@")'
$ws.Range("B$row").Formula = '=TEXT("20:18:06","@")'
$ws.Range("C$row").Formula = '=TEXT("Sunday","@")'
$ws.Range("D$row").Formula = '=TEXT("05","@")'

$textRange = $ws.Range("A${row}:D${row}")
$textRange.Copy()
$textRange.PasteSpecial(-4163)  # xlPasteValues

# --- Numeric columns (E-T) ---------------------------------------------
$ws.Range("E$row").Value = 125920
$ws.Range("F$row").Value = 141901
$ws.Range("G$row").Value = 166455
$ws.Range("H$row").Value = 157842
$ws.Range("I$row").Value = -1
$ws.Range("J$row").Value = 142077
$ws.Range("K$row").Value = -1
$ws.Range("L$row").Value = -1
$ws.Range("M$row").Value = 191157
$ws.Range("N$row").Value = 115450
$ws.Range("O$row").Value = 44805
$ws.Range("P$row").Value = 28253
$ws.Range("Q$row").Value = 63319
$ws.Range("R$row").Value = -1
$ws.Range("S$row").Value = 40151
$ws.Range("T$row").Value = -1
